$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# BLEU score (row 11)
$ws.Range("B11").Value = 0.1585103557999228

# Code BLEU (row 12)
$ws.Range("B12").Value = 0.3743823223407795
$ws.Range("C12").Value = "{'codebleu': 0.3743823223407795, 'ngram_match_score': 0.15851035579992281, 'weighted_ngram_match_score': 0.21964515665517176, 'syntax_match_score': 0.547945205479452, 'dataflow_match_score': 0.5714285714285714}"

# Embeddings and Cosine similarity (row 13)
$ws.Range("B13").Value = 0.7699755295458955
